# Apply the scene-cat block-order update:
# - Swap the B1/D1 header labels (kitchens_2 <-> living_rooms_1)
# - Update the corresponding 0/1 indicator cells so the "1" marker stays
#   aligned with the same underlying category after the header swap.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Swap header labels in row 1
$ws.Range("B1").Value = "living_rooms_1"
$ws.Range("D1").Value = "kitchens_2"

# Row 3: move the 1-marker from A3 to D3
$ws.Range("A3").Value = 0
$ws.Range("D3").Value = 1

# Row 4: move the 1-marker from D4 to B4
$ws.Range("B4").Value = 1
$ws.Range("D4").Value = 0

# Row 6: move the 1-marker from B6 to A6
$ws.Range("A6").Value = 1
$ws.Range("B6").Value = 0
